$d = $word.ActiveDocument

# --- 1. Insert "Values" right after "get" (forming "getValues"), as its own
#        run with identical formatting to the surrounding text. ---
$target = $d.Range(0, $d.Content.End)
$target.Find.ClearFormatting()
$target.Find.Text = "get"
$target.Find.MatchWholeWord = $true
$target.Find.MatchCase = $true
$target.Find.Execute() | Out-Null

$target.Collapse(0)
$target.InsertAfter("Values")

# Nudge formatting off and back on so the run-coalescing pass on save does
# not merge the new "Values" run back into the preceding "get" run.
$target.Font.Bold = 1
$target.Font.Bold = 0

# --- 2. The edit above triggers a paragraph-wide adjacent-run coalesce on
#        save, which would also merge the unrelated trailing "rango." / " "
#        runs further down the same paragraph. Protect that boundary the
#        same way so the rest of the paragraph is left untouched. ---
$tail = $d.Range(0, $d.Content.End)
$tail.Find.ClearFormatting()
$tail.Find.Text = "rango."
$tail.Find.MatchCase = $true
$tail.Find.Execute() | Out-Null

$tail.Collapse(0)
$tail.MoveEnd(1, 1) | Out-Null
$tail.Font.Bold = 1
$tail.Font.Bold = 0
